# Update the division problems/answers table cell-by-cell (by position)
# so that identical old/new text values elsewhere in the table cannot
# cause ambiguous Find & Replace matches.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$mismatches = 0

# Row 1, Col 1: "988÷6=164, 4" -> "970÷2=485, 0"
$cell = $t.Cell(1,1)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "988÷6=164, 4") {
    Write-Output "WARNING: Row 1 Col 1 expected `"988÷6=164, 4`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "970÷2=485, 0"

# Row 1, Col 2: "487÷6=81, 1" -> "627÷7=89, 4"
$cell = $t.Cell(1,2)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "487÷6=81, 1") {
    Write-Output "WARNING: Row 1 Col 2 expected `"487÷6=81, 1`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "627÷7=89, 4"

# Row 1, Col 3: "126÷5=25, 1" -> "106÷9=11, 7"
$cell = $t.Cell(1,3)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "126÷5=25, 1") {
    Write-Output "WARNING: Row 1 Col 3 expected `"126÷5=25, 1`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "106÷9=11, 7"

# Row 1, Col 4: "917÷3=305, 2" -> "939÷4=234, 3"
$cell = $t.Cell(1,4)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "917÷3=305, 2") {
    Write-Output "WARNING: Row 1 Col 4 expected `"917÷3=305, 2`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "939÷4=234, 3"

# Row 1, Col 5: "519÷9=57, 6" -> "580÷3=193, 1"
$cell = $t.Cell(1,5)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "519÷9=57, 6") {
    Write-Output "WARNING: Row 1 Col 5 expected `"519÷9=57, 6`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "580÷3=193, 1"

# Row 5, Col 1: "486÷6=81, 0" -> "687÷7=98, 1"
$cell = $t.Cell(5,1)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "486÷6=81, 0") {
    Write-Output "WARNING: Row 5 Col 1 expected `"486÷6=81, 0`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "687÷7=98, 1"

# Row 5, Col 2: "683÷9=75, 8" -> "244÷3=81, 1"
$cell = $t.Cell(5,2)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "683÷9=75, 8") {
    Write-Output "WARNING: Row 5 Col 2 expected `"683÷9=75, 8`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "244÷3=81, 1"

# Row 5, Col 3: "540÷2=270, 0" -> "895÷9=99, 4"
$cell = $t.Cell(5,3)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "540÷2=270, 0") {
    Write-Output "WARNING: Row 5 Col 3 expected `"540÷2=270, 0`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "895÷9=99, 4"

# Row 5, Col 4: "353÷4=88, 1" -> "341÷6=56, 5"
$cell = $t.Cell(5,4)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "353÷4=88, 1") {
    Write-Output "WARNING: Row 5 Col 4 expected `"353÷4=88, 1`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "341÷6=56, 5"

# Row 5, Col 5: "539÷4=134, 3" -> "978÷6=163, 0"
$cell = $t.Cell(5,5)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "539÷4=134, 3") {
    Write-Output "WARNING: Row 5 Col 5 expected `"539÷4=134, 3`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "978÷6=163, 0"

# Row 9, Col 1: "628÷8=78, 4" -> "611÷4=152, 3"
$cell = $t.Cell(9,1)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "628÷8=78, 4") {
    Write-Output "WARNING: Row 9 Col 1 expected `"628÷8=78, 4`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "611÷4=152, 3"

# Row 9, Col 2: "450÷7=64, 2" -> "771÷3=257, 0"
$cell = $t.Cell(9,2)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "450÷7=64, 2") {
    Write-Output "WARNING: Row 9 Col 2 expected `"450÷7=64, 2`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "771÷3=257, 0"

# Row 9, Col 3: "464÷5=92, 4" -> "815÷3=271, 2"
$cell = $t.Cell(9,3)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "464÷5=92, 4") {
    Write-Output "WARNING: Row 9 Col 3 expected `"464÷5=92, 4`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "815÷3=271, 2"

# Row 9, Col 4: "581÷4=145, 1" -> "281÷9=31, 2"
$cell = $t.Cell(9,4)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "581÷4=145, 1") {
    Write-Output "WARNING: Row 9 Col 4 expected `"581÷4=145, 1`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "281÷9=31, 2"

# Row 9, Col 5: "942÷9=104, 6" -> "568÷9=63, 1"
$cell = $t.Cell(9,5)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "942÷9=104, 6") {
    Write-Output "WARNING: Row 9 Col 5 expected `"942÷9=104, 6`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "568÷9=63, 1"

# Row 13, Col 1: "796÷8=99, 4" -> "557÷9=61, 8"
$cell = $t.Cell(13,1)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "796÷8=99, 4") {
    Write-Output "WARNING: Row 13 Col 1 expected `"796÷8=99, 4`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "557÷9=61, 8"

# Row 13, Col 2: "156÷6=26, 0" -> "158÷7=22, 4"
$cell = $t.Cell(13,2)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "156÷6=26, 0") {
    Write-Output "WARNING: Row 13 Col 2 expected `"156÷6=26, 0`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "158÷7=22, 4"

# Row 13, Col 3: "929÷3=309, 2" -> "163÷3=54, 1"
$cell = $t.Cell(13,3)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "929÷3=309, 2") {
    Write-Output "WARNING: Row 13 Col 3 expected `"929÷3=309, 2`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "163÷3=54, 1"

# Row 13, Col 4: "235÷3=78, 1" -> "251÷7=35, 6"
$cell = $t.Cell(13,4)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "235÷3=78, 1") {
    Write-Output "WARNING: Row 13 Col 4 expected `"235÷3=78, 1`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "251÷7=35, 6"

# Row 13, Col 5: "132÷6=22, 0" -> "555÷8=69, 3"
$cell = $t.Cell(13,5)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "132÷6=22, 0") {
    Write-Output "WARNING: Row 13 Col 5 expected `"132÷6=22, 0`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "555÷8=69, 3"

# Row 17, Col 1: "908÷7=129, 5" -> "527÷3=175, 2"
$cell = $t.Cell(17,1)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "908÷7=129, 5") {
    Write-Output "WARNING: Row 17 Col 1 expected `"908÷7=129, 5`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "527÷3=175, 2"

# Row 17, Col 2: "243÷6=40, 3" -> "988÷6=164, 4"
$cell = $t.Cell(17,2)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "243÷6=40, 3") {
    Write-Output "WARNING: Row 17 Col 2 expected `"243÷6=40, 3`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "988÷6=164, 4"

# Row 17, Col 3: "409÷7=58, 3" -> "250÷4=62, 2"
$cell = $t.Cell(17,3)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "409÷7=58, 3") {
    Write-Output "WARNING: Row 17 Col 3 expected `"409÷7=58, 3`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "250÷4=62, 2"

# Row 17, Col 4: "894÷6=149, 0" -> "142÷4=35, 2"
$cell = $t.Cell(17,4)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "894÷6=149, 0") {
    Write-Output "WARNING: Row 17 Col 4 expected `"894÷6=149, 0`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "142÷4=35, 2"

# Row 17, Col 5: "149÷7=21, 2" -> "879÷6=146, 3"
$cell = $t.Cell(17,5)
$current = $cell.Range.Text.TrimEnd([char]13, [char]7)
if ($current -ne "149÷7=21, 2") {
    Write-Output "WARNING: Row 17 Col 5 expected `"149÷7=21, 2`" but found `"$current`""
    $mismatches++
}
$cell.Range.Text = "879÷6=146, 3"

Write-Output "Updated 25 cells ($mismatches mismatches)."
